# Updated cryptos list with GitHub Actions
# Refreshes the Price (column D) and Volume(1h) (column E) figures for each
# coin row, and swaps the Filecoin/Kaspa rows (39 <-> 40) to reflect their
# new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: some "Price" values are plain decimal numbers (e.g. "601.57").
# Excel's normal auto-detection would silently turn those into numeric
# cells, so for those specific cells we force the Text number format
# first, just like a user would do when typing a numeric-looking value
# that must stay text. Values with multiple dots (e.g. "69.145.16") are
# never auto-parsed as numbers, so they don't need this treatment.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "69.145.16"
$ws.Range("E2").Value = "  +0.40%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "3.742.59"
$ws.Range("E3").Value = "  +0.19%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.57"
$ws.Range("E5").Value = "  +0.06%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.35"
$ws.Range("E6").Value = "  -0.16%  "

# Row 7 - LidoStakedEther
$ws.Range("D7").Value = "3.741.02"
$ws.Range("E7").Value = "  +0.16%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.01%  "

# Row 9 - XRP
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.541"
$ws.Range("E9").Value = "  +1.48%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +3.91%  "

# Row 11 - Toncoin
$ws.Range("E11").Value = "  +0.49%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +0.30%  "

# Row 13 - Avalanche
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.01"
$ws.Range("E13").Value = "  -0.06%  "

# Row 14 - ShibaInu
$ws.Range("E14").Value = "  +1.97%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "4.366.86"
$ws.Range("E15").Value = "  +0.29%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "3.745.83"
$ws.Range("E16").Value = "  +0.35%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "69.082.47"
$ws.Range("E17").Value = "  +0.40%  "

# Row 18 - Polkadot
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.36"
$ws.Range("E18").Value = "  +1.34%  "

# Row 19 - Chainlink
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.33"
$ws.Range("E19").Value = "  +0.38%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  -1.61%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.15"
$ws.Range("E21").Value = "  +9.47%  "

# Row 22 - BitcoinCash
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "492.96"
$ws.Range("E22").Value = "  -0.84%  "

# Row 23 - Polygon
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.727"
$ws.Range("E23").Value = "  +0.57%  "

# Row 24 - PEPE
$ws.Range("E24").Value = "  +8.24%  "

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.96"
$ws.Range("E25").Value = "  -0.25%  "

# Row 26 - Fetch.AI
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  -0.70%  "

# Row 27 - InternetComputer(DFINITY)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.25"
$ws.Range("E27").Value = "  -1.01%  "

# Row 28 - RenderToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.07"
$ws.Range("E28").Value = "  -0.56%  "

# Row 30 - PancakeSwap
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.97"
$ws.Range("E30").Value = "  +1.17%  "

# Row 31 - NEARProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.18"
$ws.Range("E31").Value = "  +2.77%  "

# Row 32 - ImmutableX
$ws.Range("E32").Value = "  +0.19%  "

# Row 33 - EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.50"
$ws.Range("E33").Value = "  -0.67%  "

# Row 34 - WrappedeETH
$ws.Range("D34").Value = "3.890.22"
$ws.Range("E34").Value = "  +0.25%  "

# Row 35 - RenzoRestakedETH
$ws.Range("D35").Value = "3.677.52"
$ws.Range("E35").Value = "  +0.41%  "

# Row 36 - Hedera
$ws.Range("E36").Value = "  -0.08%  "

# Row 37 - FirstDigitalUSD
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  +0.02%  "

# Row 38 - Mantle
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.01"
$ws.Range("E38").Value = "  +0.03%  "

# Row 39 - was Kaspa, now Filecoin (rows 39/40 swapped order)
$ws.Range("B39").Value = "Filecoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.94"
$ws.Range("E39").Value = "  +2.02%  "

# Row 40 - was Filecoin, now Kaspa
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.139"
$ws.Range("E40").Value = "  +4.62%  "

# Row 41 - TheGraph
$ws.Range("E41").Value = "  -0.25%  "

# Row 42 - dogwifhat
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.04"
$ws.Range("E42").Value = "  +6.25%  "

# Row 43 - OKB
$ws.Range("E43").Value = "  -0.31%  "

# Row 44 - Stacks
$ws.Range("E44").Value = "  +0.27%  "

# Row 45 - Bittensor
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "424.31"
$ws.Range("E45").Value = "  -2.62%  "

# Row 46 - Cosmos
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.45"
$ws.Range("E46").Value = "  +0.55%  "

# Row 48 - Arweave
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.07"
$ws.Range("E48").Value = "  -1.46%  "

# Row 49 - Monero
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.14"
$ws.Range("E49").Value = "  +0.39%  "

# Row 50 - Maker
$ws.Range("D50").Value = "2.781.05"
$ws.Range("E50").Value = "  +1.41%  "

# Row 51 - VeChain
$ws.Range("E51").Value = "  +0.33%  "
